# This workbook is a weekly price log for "Frutilla" (strawberry) sold at
# "Vega Central Mapocho de Santiago". A new week of data (4 quality-grade
# rows, dated 44474) is inserted right above the existing block that starts
# at row 446, pushing all the following rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at the top of the insertion point; this shifts every
# row currently at 446 and below down by 4 (446 -> 450, ..., 523 -> 527).
$ws.Range("A446:A449").EntireRow.Insert()

# Common values shared by every data row in this sheet.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$kgUnidad    = 7

function Set-FrutillaRow($RowNum, $Fecha, $Calidad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Origen, $PrecioKg) {
    $ws.Cells.Item($RowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($RowNum, 2).Value  = $mercado
    $ws.Cells.Item($RowNum, 3).Value  = $region
    $ws.Cells.Item($RowNum, 4).Value  = $Fecha
    $ws.Cells.Item($RowNum, 5).Value  = $codreg
    $ws.Cells.Item($RowNum, 6).Value  = $tipo
    $ws.Cells.Item($RowNum, 7).Value  = $productoId
    $ws.Cells.Item($RowNum, 8).Value  = $producto
    $ws.Cells.Item($RowNum, 9).Value  = $categoriaId
    $ws.Cells.Item($RowNum, 10).Value = $categoria
    $ws.Cells.Item($RowNum, 11).Value = $variedad
    $ws.Cells.Item($RowNum, 12).Value = $Calidad
    $ws.Cells.Item($RowNum, 13).Value = $Volumen
    $ws.Cells.Item($RowNum, 14).Value = $PrecioMinimo
    $ws.Cells.Item($RowNum, 15).Value = $PrecioMaximo
    $ws.Cells.Item($RowNum, 16).Value = $PrecioPromedio
    $ws.Cells.Item($RowNum, 17).Value = $unidad
    $ws.Cells.Item($RowNum, 18).Value = $Origen
    $ws.Cells.Item($RowNum, 19).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value = $kgUnidad
}

Set-FrutillaRow 446 44474 "Especial" 590 11000 12000 11475 "Provincia de Melipilla" 1639
Set-FrutillaRow 447 44474 "Primera"  650 9000  10000 9462  "Provincia de Melipilla" 1352
Set-FrutillaRow 448 44474 "Segunda"  530 7000  8000  7528  "Provincia de Melipilla" 1075
Set-FrutillaRow 449 44474 "Tercera"  530 3000  4000  3528  "Provincia de Melipilla" 504

Write-Host "Inserted 4 new Frutilla rows (446-449) for date 44474"
